$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7851601243019104
$ws.Range("B1").Value = 3.580848217010498
$ws.Range("C1").Value = 3.487923622131348
$ws.Range("D1").Value = 2.908017873764038
$ws.Range("E1").Value = 1.78067684173584
